# Add 2022-Q4 data:
#  - "总计" sheet gets a new row 2 for 2022-Q4 (existing 2022-Q3 row shifts to row 3)
#  - the current "2022-Q3" sheet (per-fund holdings) is duplicated so the
#    original data survives under the name "2022-Q3"
#  - the original sheet is renamed to "2022-Q4" and repopulated with the new
#    quarter's fund data

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Item(2)

# ---- 总计: insert the 2022-Q4 row above the existing 2022-Q3 row ----
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 5
$wsTotal.Cells.Item(2, 4).Value = 0.14

# ---- duplicate the current per-fund sheet so 2022-Q3's data is preserved ----
$wsQ3.Copy($null, $wsQ3)
$wsOldQ3 = $wb.Worksheets.Item(3)
$wsQ3.Name = "2022-Q4"
$wsOldQ3.Name = "2022-Q3"

# ---- repopulate the (renamed) sheet with the new 2022-Q4 fund data ----
$wsQ4 = $wsQ3

# helper cells on 总计 carry the two formats we need to reuse:
#   B1 -> bold/centered "header" style
#   D2 -> the plain/default style
$headerFmt = $wsTotal.Range("B1")
$plainFmt = $wsTotal.Range("D2")

# stamp the header style first (format-only paste), then write the header text
$headerFmt.Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsQ4.Cells.Item(1, 2).Value = "基金代码"
$wsQ4.Cells.Item(1, 3).Value = "基金名称"
$wsQ4.Cells.Item(1, 4).Value = "基金规模"
$wsQ4.Cells.Item(1, 5).Value = "股票总仓位"
$wsQ4.Cells.Item(1, 6).Value = "仓位占比"
$wsQ4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$wsQ4.Cells.Item(1, 8).Value = "仓位排名"

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2
$wsQ4.Cells.Item(2, 1).Value = 0
Set-TextValue $wsQ4.Cells.Item(2, 2) "015071"
$wsQ4.Cells.Item(2, 3).Value = "鑫元专精特新混合A"
Set-TextValue $wsQ4.Cells.Item(2, 4) "2.46"
Set-TextValue $wsQ4.Cells.Item(2, 5) "70.69"
Set-TextValue $wsQ4.Cells.Item(2, 6) "3.70"
Set-TextValue $wsQ4.Cells.Item(2, 7) "0.0910"
$wsQ4.Cells.Item(2, 8).Value = 6

# Row 3
$wsQ4.Cells.Item(3, 1).Value = 1
Set-TextValue $wsQ4.Cells.Item(3, 2) "005075"
$wsQ4.Cells.Item(3, 3).Value = "富国研究量化精选混合"
Set-TextValue $wsQ4.Cells.Item(3, 4) "2.48"
Set-TextValue $wsQ4.Cells.Item(3, 5) "90.71"
Set-TextValue $wsQ4.Cells.Item(3, 6) "1.68"
Set-TextValue $wsQ4.Cells.Item(3, 7) "0.0417"
$wsQ4.Cells.Item(3, 8).Value = 2

# Row 4
$wsQ4.Cells.Item(4, 1).Value = 2
Set-TextValue $wsQ4.Cells.Item(4, 2) "015072"
$wsQ4.Cells.Item(4, 3).Value = "鑫元专精特新混合C"
Set-TextValue $wsQ4.Cells.Item(4, 4) "0.25"
Set-TextValue $wsQ4.Cells.Item(4, 5) "70.69"
Set-TextValue $wsQ4.Cells.Item(4, 6) "3.70"
Set-TextValue $wsQ4.Cells.Item(4, 7) "0.0092"
$wsQ4.Cells.Item(4, 8).Value = 6

# Row 5
$wsQ4.Cells.Item(5, 1).Value = 3
Set-TextValue $wsQ4.Cells.Item(5, 2) "011494"
$wsQ4.Cells.Item(5, 3).Value = "华泰紫金丰和偏债混合发起A"
Set-TextValue $wsQ4.Cells.Item(5, 4) "0.12"
Set-TextValue $wsQ4.Cells.Item(5, 5) "36.70"
Set-TextValue $wsQ4.Cells.Item(5, 6) "1.37"
Set-TextValue $wsQ4.Cells.Item(5, 7) "0.0016"
$wsQ4.Cells.Item(5, 8).Value = 4

# Row 6
$wsQ4.Cells.Item(6, 1).Value = 4
Set-TextValue $wsQ4.Cells.Item(6, 2) "011495"
$wsQ4.Cells.Item(6, 3).Value = "华泰紫金丰和偏债混合发起C"
Set-TextValue $wsQ4.Cells.Item(6, 4) "0.03"
Set-TextValue $wsQ4.Cells.Item(6, 5) "36.70"
Set-TextValue $wsQ4.Cells.Item(6, 6) "1.37"
Set-TextValue $wsQ4.Cells.Item(6, 7) "0.0004"
$wsQ4.Cells.Item(6, 8).Value = 4

# restore the plain/default format on the text-forced cells (B:G except
# the "A" index column, which keeps the bold/centered style below) -
# format-only paste so the values just written above are preserved
$plainFmt.Copy()
$wsQ4.Range("B2:G6").PasteSpecial(-4122)

# column A (row index) keeps the bold/centered header-like style too
$wsHeaderA = $wsTotal.Range("A2")
$wsHeaderA.Copy()
$wsQ4.Range("A2:A6").PasteSpecial(-4122)
$wsQ4.Cells.Item(2, 1).Value = 0
$wsQ4.Cells.Item(3, 1).Value = 1
$wsQ4.Cells.Item(4, 1).Value = 2
$wsQ4.Cells.Item(5, 1).Value = 3
$wsQ4.Cells.Item(6, 1).Value = 4
